$wb = $excel.ActiveWorkbook

# Delete the first sheet "papa"
$papa = $wb.Worksheets.Item("papa")
$papa.Delete()

# Rename "hola" to "Algo"
$hola = $wb.Worksheets.Item("hola")
$hola.Name = "Algo"

# Update PINI sheet's second row values
$pini = $wb.Worksheets.Item("PINI")
$pini.Range("A2").Value = "Sabana"
$pini.Range("B2").NumberFormat = "@"
$pini.Range("B2").Value = "2"
